# Applies the edit described in the commit "Match All Report SpreadSheets"
# to the Tax Summary report workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the main report sheet (title-case it).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Tax Summary report result")
$ws1.Name = "Tax Summary Report Result"

$ws2 = $wb.Worksheets.Item("Raw Data")

$tableName = "TTaxSummaryReport_IgnoreDates_false_ReportType_Summary_DateFrom__222022_12_31_22"

# ---------------------------------------------------------------------
# 2. Rebuild the report rows on the main sheet as live formulas that
#    pull from the "Raw Data" table, instead of the old hard-coded
#    values (which were also shifted one column to the right of their
#    header, i.e. column B/.../G instead of A/.../F).
# ---------------------------------------------------------------------
$rows = @(
    @{ Row = 2; Code = "CAG" },
    @{ Row = 3; Code = "FRE" },
    @{ Row = 4; Code = "GST" },
    @{ Row = 5; Code = "NCG" },
    @{ Row = 6; Code = "WC" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws1.Range("A$row").Formula = "=$tableName[[#This Row],[T.TaxCode]]"
    $ws1.Range("B$row").Formula = "=$tableName[[#This Row],[T.INPUT_AmountEx]]"
    $ws1.Range("C$row").Formula = "=$tableName[[#This Row],[T.OUTPUT_AmountEx]]"
    $ws1.Range("D$row").Formula = "=$tableName[[#This Row],[T.INPUT_AmountInc]]"
    $ws1.Range("E$row").Formula = "=$tableName[[#This Row],[T.TaxRate]]"
    $ws1.Range("F$row").Formula = "=$tableName[[#This Row],[T.TotalTax]]"

    # Currency formatting for the amount columns, percentage for the rate.
    $ws1.Range("B$row").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
    $ws1.Range("C$row").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
    $ws1.Range("D$row").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
    $ws1.Range("F$row").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
    $ws1.Range("E$row").NumberFormat = "0.00%"
}

# Column G is no longer used now that the data re-aligned under its
# real headers (A-F) -- clear it entirely so the sheet dimension
# shrinks back down to F7.
$ws1.Range("G2:G6").Clear()

# Grand-total row: turn the static totals into SUM formulas.
$ws1.Range("D7").Formula = "=SUM(D2:D6)"
$ws1.Range("F7").Formula = "=SUM(F2:F6)"

# ---------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping to match the saved view state.
# ---------------------------------------------------------------------
$ws1.Range("F7").Select()
$ws2.Range("J5").Select()
$ws1.Activate()

# ---------------------------------------------------------------------
# 4. Table style refresh on the "Raw Data" table.
# ---------------------------------------------------------------------
$table = $ws2.ListObjects.Item(1)
$table.TableStyle = "TableStyleMedium7"

$wb.RefreshAll()
